$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.241.94'
$ws.Range("E2").Value = '  +0.35%  '

$ws.Range("D3").Value = '3.802.70'
$ws.Range("E3").Value = '  -0.76%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '707.85'
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.17'
$ws.Range("E6").Value = '  -0.90%  '

$ws.Range("D7").Value = '3.797.90'
$ws.Range("E7").Value = '  -0.85%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("E10").Value = '  -1.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.50'
$ws.Range("E11").Value = '  +2.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.480'
$ws.Range("E12").Value = '  +4.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  -1.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.23'
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("D15").Value = '4.440.19'
$ws.Range("E15").Value = '  -0.86%  '

$ws.Range("D16").Value = '3.786.96'
$ws.Range("E16").Value = '  -0.84%  '

$ws.Range("D17").Value = '71.233.70'
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.18'
$ws.Range("E18").Value = '  -0.41%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.50'
$ws.Range("E19").Value = '  +0.94%  '

$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.115'
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '516.55'
$ws.Range("E21").Value = '  +4.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.43'
$ws.Range("E22").Value = '  -2.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.715'
$ws.Range("E23").Value = '  -0.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.08'
$ws.Range("E24").Value = '  -1.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000141'
$ws.Range("E25").Value = '  -3.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.66'
$ws.Range("E26").Value = '  +4.21%  '

$ws.Range("D27").Value = '3.943.63'
$ws.Range("E27").Value = '  -1.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  -2.66%  '

$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("E30").Value = '  -3.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  -5.07%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.37'
$ws.Range("E32").Value = '  -1.63%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.24'
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.12'
$ws.Range("E34").Value = '  -0.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.173'
$ws.Range("E35").Value = '  -2.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.26'
$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").Value = '3.764.73'
$ws.Range("E38").Value = '  -0.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.46'
$ws.Range("E39").Value = '  +7.57%  '

$ws.Range("E40").Value = '  -1.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.44'
$ws.Range("E41").Value = '  +5.45%  '

$ws.Range("E42").Value = '  -1.93%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.25'
$ws.Range("E43").Value = '  -3.00%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '168.33'
$ws.Range("E46").Value = '  +2.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '50.06'
$ws.Range("E47").Value = '  +2.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000307'
$ws.Range("E48").Value = '  -1.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '427.61'
$ws.Range("E49").Value = '  +3.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.40'
$ws.Range("E50").Value = '  +0.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.60'
$ws.Range("E51").Value = '  -0.09%  '
